$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.57077
$ws.Range("H2").Value = 4.71231
$ws.Range("I2").Value = 0.02582502173444737
$ws.Range("J2").Value = 0.02582502173444737
$ws.Range("Q2").Value = 0.09800819415000002
$ws.Range("R2").Value = 0.88207374735
$ws.Range("S2").Value = 0.02582502173444737
$ws.Range("T2").Value = 0.02582502173444737

$ws.Range("I3").Value = 0.934831682683009
$ws.Range("J3").Value = 0.934831682683009
$ws.Range("S3").Value = 0.934831682683009
$ws.Range("T3").Value = 0.934831682683009

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9273763333333335
$ws.Range("H4").Value = 2.782129
$ws.Range("I4").Value = 0.01524698967025436
$ws.Range("J4").Value = 0.01524698967025436
$ws.Range("Q4").Value = 0.05786364631833334
$ws.Range("R4").Value = 0.520772816865
$ws.Range("S4").Value = 0.01524698967025436
$ws.Range("T4").Value = 0.01524698967025436

$ws.Range("G5").Value = 0.7810079999999999
$ws.Range("H5").Value = 2.343024
$ws.Range("I5").Value = 0.0128405486320577
$ws.Range("J5").Value = 0.0128405486320577
$ws.Range("Q5").Value = 0.04873099415999999
$ws.Range("R5").Value = 0.4385789474399999
$ws.Range("S5").Value = 0.0128405486320577
$ws.Range("T5").Value = 0.0128405486320577

$ws.Range("G6").Value = 0.6846153333333334
$ws.Range("H6").Value = 2.053846
$ws.Range("I6").Value = 0.01125575728023152
$ws.Range("J6").Value = 0.01125575728023152
$ws.Range("Q6").Value = 0.04271657372333333
$ws.Range("R6").Value = 0.38444916351
$ws.Range("S6").Value = 0.01125575728023152
$ws.Range("T6").Value = 0.01125575728023152
